$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3701067615658363
$ws1.Range("C2").Value = 0.07329842931937172
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1365853658536585
$ws1.Range("F2").Value = 0.2834008097165992
$ws1.Range("G2").Value = 0.6728280961182994
$ws1.Range("H2").Value = 0.8257758159443552
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 354
$ws1.Range("K2").Value = 180
$ws1.Range("L2").Value = 0

# ---- Sheet 2: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 - label "0"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.3370786516853932
$ws2.Range("D2").Value = 0.5042016806722689

# row 3 - label "1"
$ws2.Range("B3").Value = 0.07329842931937172
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1365853658536585

# row 4 - label "accuracy"
$ws2.Range("B4").Value = 0.3701067615658363
$ws2.Range("C4").Value = 0.3701067615658363
$ws2.Range("D4").Value = 0.3701067615658363
$ws2.Range("E4").Value = 0.3701067615658363

# row 5 - label "macro avg"
$ws2.Range("B5").Value = 0.5366492146596858
$ws2.Range("C5").Value = 0.6685393258426966
$ws2.Range("D5").Value = 0.3203935232629637

# row 6 - label "weighted avg"
$ws2.Range("B6").Value = 0.9538298149838833
$ws2.Range("C6").Value = 0.3701067615658363
$ws2.Range("D6").Value = 0.485886277087

# ---- Sheet 3: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 - Actual 0
$ws3.Range("B2").Value = 180
$ws3.Range("C2").Value = 354

# row 3 - Actual 1
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
